# 141: 30/12 18:35 update - LP1912 / LP1912-215 / 6203-6173 sheets
$wb = $excel.ActiveWorkbook

function Fill-Rows($ws, [string]$rangeAddr, $rowsData) {
    $nRows = $rowsData.Count
    $nCols = 7
    $arr = New-Object 'object[,]' $nRows,$nCols
    for ($i = 0; $i -lt $nRows; $i++) {
        $r = $rowsData[$i]
        for ($j = 0; $j -lt $nCols; $j++) {
            $arr[$i,$j] = $r[$j]
        }
    }
    $ws.Range($rangeAddr).Value2 = $arr
}

# ---------------------------------------------------------------------
# Sheet 1: LP1912  (cols: A(blank) B=Hora_Scrap C=Hora_Llegada D=Linea E=Minutos F=Parada G=Fecha)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value2 = "Última actualización: 30/12/2025 15:35:22"
$ws1.Range("A3").Value2 = "Total filas: 338"

$s1rows = @(
    @("", "15:35:11", "15:38", "23_HERNANDEZ", 3, "LP1912", "30/12/2025"),
    @("", "15:35:11", "15:46", "14_ABASTO", 11, "LP1912", "30/12/2025"),
    @("", "15:35:11", "15:54", "11_ETCHEVERRY", 19, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:01", "10_OLMOS", 26, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:03", "16_SANTA ANA", 28, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:05", "23_HERNANDEZ", 30, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:11", "16_SANTA ANA", 36, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:20", "215C_EL PATO", 45, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:21", "26_HERNANDEZ", 46, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:29", "10_OLMOS", 54, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:37", "11_ETCHEVERRY", 62, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:43", "16_P MOR-SANTA ANA", 68, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:45", "14_ABASTO", 70, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:45", "23_HERNANDEZ", 70, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:48", "15_ABASTO", 73, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:56", "17_179 Y 38", 81, "LP1912", "30/12/2025"),
    @("", "15:35:11", "16:57", "10_OLMOS", 82, "LP1912", "30/12/2025"),
    @("", "15:35:11", "17:05", "215A_EL PATO", 90, "LP1912", "30/12/2025"),
    @("", "15:35:11", "17:11", "11_ETCHEVERRY", 96, "LP1912", "30/12/2025")
)
Fill-Rows $ws1 "A321:G339" $s1rows

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215 (cols: A(blank) B=Fecha C=Hora_Scrap D=Hora_Llegada E=Linea F=Minutos G=Parada)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value2 = "Última actualización: 30/12/2025 15:35:22"
$ws2.Range("A3").Value2 = "Total filas: 26"

$s2rows = @(
    @("", "30/12/2025", "15:35:11", "16:20", "215C_EL PATO", 45, "LP1912"),
    @("", "30/12/2025", "15:35:11", "17:05", "215A_EL PATO", 90, "LP1912")
)
Fill-Rows $ws2 "A26:G27" $s2rows

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173 (cols: A(blank) B=Fecha C=Hora_Scrap D=Hora_Llegada E=Linea F=Minutos G=Parada)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value2 = "Última actualización: 30/12/2025 15:35:22"
$ws3.Range("A3").Value2 = "Total filas: 47"

$s3rows = @(
    @("", "30/12/2025", "15:35:17", "16:14", "215C_LA PLATA", 39, "L6203"),
    @("", "30/12/2025", "15:35:22", "16:53", "215B_LP-P MOR-40 Y 115", 78, "L6173")
)
Fill-Rows $ws3 "A47:G48" $s3rows
